# Insert a new data row above the current row 129 (pushing the existing
# rows 129-206 down to 130-207) and populate it with a new price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 129..206 down to 130..207, creating a blank row 129.
$ws.Rows("129:129").Insert()

# Fill the new row 129 with the new record's data.
$ws.Range("A129").Value = 7
$ws.Range("B129").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C129").Value = "Ñuble"
$ws.Range("D129").Value = 44596
$ws.Range("E129").Value = 16
$ws.Range("F129").Value = 100112009
$ws.Range("G129").Value = "Acelga"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 100
$ws.Range("K129").Value = 400
$ws.Range("L129").Value = 450
$ws.Range("M129").Value = 425
$ws.Range("N129").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O129").Value = "Provincia de Diguillín"
$ws.Range("P129").Value = 425
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"
